$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "51.686.36"
$ws.Range("E2").Value = "  +5.69%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.759.15"
$ws.Range("E3").Value = "  +5.30%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  +0.13%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "117.93"
$ws.Range("E5").Value = "  +6.94%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "332.72"
$ws.Range("E6").Value = "  +3.96%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.535"
$ws.Range("E7").Value = "  +3.14%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  +0.01%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.578"
$ws.Range("E9").Value = "  +7.24%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.76"
$ws.Range("E10").Value = "  +6.06%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.24"
$ws.Range("E11").Value = "  +2.50%  "

$ws.Range("E12").Value = "  +3.24%  "

$ws.Range("E13").Value = "  +3.31%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.67"
$ws.Range("E14").Value = "  +6.78%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.190.51"
$ws.Range("E15").Value = "  +5.16%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.755.62"
$ws.Range("E16").Value = "  +5.51%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.886"
$ws.Range("E17").Value = "  +4.05%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "51.572.83"
$ws.Range("E18").Value = "  +5.64%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.73"
$ws.Range("E19").Value = "  +7.30%  "

$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.88"
$ws.Range("E20").Value = "  +3.89%  "

$ws.Range("B21").Value = "ImmutableX"
$ws.Range("C21").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.99"
$ws.Range("E21").Value = "  +4.31%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0₃0965"
$ws.Range("E22").Value = "  +2.91%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "276.98"
$ws.Range("E23").Value = "  +2.86%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "70.13"
$ws.Range("E24").Value = "  +0.35%  "

$ws.Range("E25").Value = "  +6.02%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.93"
$ws.Range("E26").Value = "  +2.95%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "4.15"
$ws.Range("E27").Value = "  +0.82%  "

$ws.Range("E28").Value = "  +0.04%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.30"
$ws.Range("E29").Value = "  +2.31%  "

$ws.Range("E30").Value = "  +0.13%  "

$ws.Range("E31").Value = "  +2.84%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "35.78"
$ws.Range("E32").Value = "  +1.25%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "50.51"
$ws.Range("E33").Value = "  +2.65%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.62"
$ws.Range("E34").Value = "  +4.28%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0830"
$ws.Range("E35").Value = "  +5.30%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "19.39"
$ws.Range("E36").Value = "  +2.00%  "

$ws.Range("E37").Value = "  +0.04%  "

$ws.Range("E38").Value = "  +4.97%  "

$ws.Range("E39").Value = "  +2.05%  "

$ws.Range("E40").Value = "  +5.38%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "130.23"
$ws.Range("E41").Value = "  +4.88%  "

$ws.Range("E42").Value = "  +5.59%  "

$ws.Range("E43").Value = "  +3.43%  "

$ws.Range("E44").Value = "  +10.36%  "

$ws.Range("E45").Value = "  +4.07%  "

$ws.Range("E46").Value = "  +12.75%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.121.02"
$ws.Range("E47").Value = "  +2.09%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.37"
$ws.Range("E48").Value = "  +5.03%  "

$ws.Range("E49").Value = "  +3.54%  "

$ws.Range("E50").Value = "  +8.33%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "9.00"
$ws.Range("E51").Value = "  +1.87%  "
